# The deck's one and only slide master currently uses the "Integral"
# (Red Violet) theme (ppt/theme/theme1.xml). A stock default "Office
# Theme" colour/font/format scheme already ships in the package, unused,
# as theme2.xml (wired up only to the notes master). The author picked a
# new Design from PowerPoint's Design tab, swapping the active theme's
# palette over to that default "Office Theme" look.
#
# The font scheme and the shape/effect format scheme are byte-for-byte
# identical between the two themes, so the only visible difference is
# the 12-slot theme colour scheme. We repaint it through the live
# ThemeColorScheme on the presentation's slide master/theme - the
# supported, persisted way to change the active theme's palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as OLE RGB() values (0x00BBGGRR) to match the hex codes from
# the target theme XML:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
$colors.Item(1).RGB  = 0          # dk1
$colors.Item(2).RGB  = 16777215   # lt1
$colors.Item(3).RGB  = 6968388    # dk2
$colors.Item(4).RGB  = 15132391   # lt2
$colors.Item(5).RGB  = 13998939   # accent1
$colors.Item(6).RGB  = 3243501    # accent2
$colors.Item(7).RGB  = 10855845   # accent3
$colors.Item(8).RGB  = 49407      # accent4
$colors.Item(9).RGB  = 12874308   # accent5
$colors.Item(10).RGB = 4697456    # accent6
$colors.Item(11).RGB = 12673797   # hlink
$colors.Item(12).RGB = 7491477    # folHlink

# Best-effort: line up the display names with the Office Theme as well.
# These are read-only in this host (and don't touch any other parts),
# so don't let a failure here drop the colour-scheme repaint above.
try { $colors.Name = "Office" } catch { }
try { $theme.Name = "Office Theme" } catch { }
